# Update countries & provincias Spain
# Applies:
#  - Barein overtakes Nueva Zelanda in the "Casos totales" ranking (rows 65/66 swap)
#  - Zimbabue overtakes several countries and moves up to just above Mongolia
#    (rows 177-186 shift down by one, Zimbabue's fresh numbers land in row 177)
#  - Noruega's daily numbers are refreshed (row 32)
#  - The "Datos actualizados" timestamp footer is refreshed (A1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row {
    param([int]$Row, [string]$Pais, [double]$CasosTotales, [double]$NuevosCasos,
          [double]$CasosActivos, [double]$Recuperados, [double]$CasosCriticos,
          [double]$MuertesHoy, [double]$Muertes)
    $ws.Cells.Item($Row, 1).Value = $Pais
    $ws.Cells.Item($Row, 2).Value = $CasosTotales
    $ws.Cells.Item($Row, 3).Value = $NuevosCasos
    $ws.Cells.Item($Row, 4).Value = $CasosActivos
    $ws.Cells.Item($Row, 5).Value = $Recuperados
    $ws.Cells.Item($Row, 6).Value = $CasosCriticos
    $ws.Cells.Item($Row, 7).Value = $MuertesHoy
    $ws.Cells.Item($Row, 8).Value = $Muertes
}

# Timestamp footer
$ws.Range("A1").Value = "Datos actualizados a 13 de Abril de 2020 a las 22:52"

# Noruega refreshed totals (ranking position unchanged)
Set-Row 32 "Noruega" 6565 40 32 6399 59 6 134

# Barein jumps ahead of Nueva Zelanda
Set-Row 65 "Barein" 1361 225 591 764 4 0 6
Set-Row 66 "Nueva Zelanda" 1349 19 546 798 4 1 5

# Zimbabue jumps ahead of Mongolia, Fiyi, Malaui, Namibia, Dominica,
# Santa Lucia, Suazilandia, Granada and Nepal (everybody else shifts down one spot)
Set-Row 177 "Zimbabue" 17 3 0 14 0 0 3
Set-Row 178 "Mongolia" 17 1 4 13 0 0 0
Set-Row 179 "Fiyi" 16 0 0 16 0 0 0
Set-Row 180 "Malaui" 16 3 0 14 1 0 2
Set-Row 181 "Namibia" 16 0 3 13 0 0 0
Set-Row 182 "Dominica" 16 0 8 8 0 0 0
Set-Row 183 "Santa Lucia" 15 0 4 11 0 0 0
Set-Row 184 "Suazilandia" 15 1 7 8 0 0 0
Set-Row 185 "Granada" 14 0 0 14 2 0 0
Set-Row 186 "Nepal" 14 2 1 13 0 0 0
